$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39 (pushes existing rows 39-69 down to 40-70)
$ws.Rows.Item(39).Insert()

# Only column D carries a non-default style (date format) - copy just that cell's number format
$ws.Cells.Item(39, 4).NumberFormat = $ws.Cells.Item(38, 4).NumberFormat

# Fill in the new row's values
$ws.Cells.Item(39, 1).Value = 10
$ws.Cells.Item(39, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(39, 3).Value = "La Araucanía"
$ws.Cells.Item(39, 4).Value = 44741
$ws.Cells.Item(39, 5).Value = 9
$ws.Cells.Item(39, 6).Value = 100112026
$ws.Cells.Item(39, 7).Value = "Haba"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 40
$ws.Cells.Item(39, 11).Value = 20000
$ws.Cells.Item(39, 12).Value = 20000
$ws.Cells.Item(39, 13).Value = 20000
$ws.Cells.Item(39, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(39, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(39, 16).Value = 800
$ws.Cells.Item(39, 17).Value = 25
$ws.Cells.Item(39, 18).Value = "Hortaliza"
